$wb = $excel.ActiveWorkbook

# --- Add the "metadata" sheet right after "data" ---
$dataSheet = $wb.Worksheets.Item("data")
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

# Match the "data" sheet's page layout (Excel defaults differ for a brand-new sheet)
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# --- Header row ---
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# Copy the header formatting (bold, bordered, centered) from the "data" sheet's header
$dataSheet.Range("B1:F1").Copy()
$ws.Range("B1:F1").PasteSpecial(-4122)
$dataSheet.Range("B1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

# --- Data row ---
$ws.Range("A2").Value = 0
$dataSheet.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("B2").Value = "Severe microcephaly"
$ws.Range("C2").Value = 162
# Force "2.256" to stay literal text (not be parsed as a number), then drop back to the
# default (unstyled) look so it matches a plain, un-styled text cell.
$ws.Range("D2").Value = "'2.256"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "2021-10-01T14:36:26.374873Z"
$ws.Range("F2").Value = "2021-10-05 14:22:40.895790"
$ws.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/162/?format=json"

$excel.CutCopyMode = $false

# Keep "data" as the active/selected sheet (matches activeTab=0 in the target workbook)
$dataSheet.Activate()

# --- Refresh the "time_taken" (column F) timestamps on the "data" sheet ---
$ds = $dataSheet
$ds.Range("F2").Value = "2021-10-05 14:22:40.899095"
$ds.Range("F3").Value = "2021-10-05 14:22:40.899103"
$ds.Range("F4").Value = "2021-10-05 14:22:40.899106"
$ds.Range("F5").Value = "2021-10-05 14:22:40.899109"
$ds.Range("F6").Value = "2021-10-05 14:22:40.899112"
$ds.Range("F7").Value = "2021-10-05 14:22:40.899115"
$ds.Range("F8").Value = "2021-10-05 14:22:40.899118"
$ds.Range("F9").Value = "2021-10-05 14:22:40.899120"
$ds.Range("F10").Value = "2021-10-05 14:22:40.899123"
$ds.Range("F11").Value = "2021-10-05 14:22:40.899126"
$ds.Range("F12").Value = "2021-10-05 14:22:40.899128"
$ds.Range("F13").Value = "2021-10-05 14:22:40.899131"
$ds.Range("F14").Value = "2021-10-05 14:22:40.899133"
$ds.Range("F15").Value = "2021-10-05 14:22:40.899136"
$ds.Range("F16").Value = "2021-10-05 14:22:40.899138"
$ds.Range("F17").Value = "2021-10-05 14:22:40.899141"
$ds.Range("F18").Value = "2021-10-05 14:22:40.899144"
$ds.Range("F19").Value = "2021-10-05 14:22:40.899147"
$ds.Range("F20").Value = "2021-10-05 14:22:40.899149"
$ds.Range("F21").Value = "2021-10-05 14:22:40.899152"
$ds.Range("F22").Value = "2021-10-05 14:22:40.899154"
$ds.Range("F23").Value = "2021-10-05 14:22:40.899157"
$ds.Range("F24").Value = "2021-10-05 14:22:40.899159"
$ds.Range("F25").Value = "2021-10-05 14:22:40.899162"
$ds.Range("F26").Value = "2021-10-05 14:22:40.899165"
$ds.Range("F27").Value = "2021-10-05 14:22:40.899168"
$ds.Range("F28").Value = "2021-10-05 14:22:40.899170"
$ds.Range("F29").Value = "2021-10-05 14:22:40.899173"
$ds.Range("F30").Value = "2021-10-05 14:22:40.899175"
$ds.Range("F31").Value = "2021-10-05 14:22:40.899178"
$ds.Range("F32").Value = "2021-10-05 14:22:40.899181"
$ds.Range("F33").Value = "2021-10-05 14:22:40.899183"
$ds.Range("F34").Value = "2021-10-05 14:22:40.899186"
$ds.Range("F35").Value = "2021-10-05 14:22:40.899189"
$ds.Range("F36").Value = "2021-10-05 14:22:40.899191"
$ds.Range("F37").Value = "2021-10-05 14:22:40.899194"
$ds.Range("F38").Value = "2021-10-05 14:22:40.899196"
$ds.Range("F39").Value = "2021-10-05 14:22:40.899199"
$ds.Range("F40").Value = "2021-10-05 14:22:40.899201"
$ds.Range("F41").Value = "2021-10-05 14:22:40.899204"
$ds.Range("F42").Value = "2021-10-05 14:22:40.899207"
$ds.Range("F43").Value = "2021-10-05 14:22:40.899209"
$ds.Range("F44").Value = "2021-10-05 14:22:40.899212"
$ds.Range("F45").Value = "2021-10-05 14:22:40.899214"
$ds.Range("F46").Value = "2021-10-05 14:22:40.899217"
$ds.Range("F47").Value = "2021-10-05 14:22:40.899219"
$ds.Range("F48").Value = "2021-10-05 14:22:40.899222"
$ds.Range("F49").Value = "2021-10-05 14:22:40.899224"
$ds.Range("F50").Value = "2021-10-05 14:22:40.899227"
$ds.Range("F51").Value = "2021-10-05 14:22:40.899229"
$ds.Range("F52").Value = "2021-10-05 14:22:40.899232"
$ds.Range("F53").Value = "2021-10-05 14:22:40.899234"
$ds.Range("F54").Value = "2021-10-05 14:22:40.899238"
$ds.Range("F55").Value = "2021-10-05 14:22:40.899240"
$ds.Range("F56").Value = "2021-10-05 14:22:40.899243"
$ds.Range("F57").Value = "2021-10-05 14:22:40.899245"
$ds.Range("F58").Value = "2021-10-05 14:22:40.899248"
$ds.Range("F59").Value = "2021-10-05 14:22:40.899250"
$ds.Range("F60").Value = "2021-10-05 14:22:40.899253"
$ds.Range("F61").Value = "2021-10-05 14:22:40.899255"
$ds.Range("F62").Value = "2021-10-05 14:22:40.899258"
$ds.Range("F63").Value = "2021-10-05 14:22:40.899261"
$ds.Range("F64").Value = "2021-10-05 14:22:40.899263"
$ds.Range("F65").Value = "2021-10-05 14:22:40.899266"
$ds.Range("F66").Value = "2021-10-05 14:22:40.899269"
$ds.Range("F67").Value = "2021-10-05 14:22:40.899272"
$ds.Range("F68").Value = "2021-10-05 14:22:40.899275"
$ds.Range("F69").Value = "2021-10-05 14:22:40.899277"
$ds.Range("F70").Value = "2021-10-05 14:22:40.899280"
$ds.Range("F71").Value = "2021-10-05 14:22:40.899282"
$ds.Range("F72").Value = "2021-10-05 14:22:40.899285"
$ds.Range("F73").Value = "2021-10-05 14:22:40.899287"
$ds.Range("F74").Value = "2021-10-05 14:22:40.899290"
$ds.Range("F75").Value = "2021-10-05 14:22:40.899292"
$ds.Range("F76").Value = "2021-10-05 14:22:40.899295"
$ds.Range("F77").Value = "2021-10-05 14:22:40.899298"
$ds.Range("F78").Value = "2021-10-05 14:22:40.899302"
$ds.Range("F79").Value = "2021-10-05 14:22:40.899305"
$ds.Range("F80").Value = "2021-10-05 14:22:40.899308"
$ds.Range("F81").Value = "2021-10-05 14:22:40.899311"
$ds.Range("F82").Value = "2021-10-05 14:22:40.899313"
$ds.Range("F83").Value = "2021-10-05 14:22:40.899316"
$ds.Range("F84").Value = "2021-10-05 14:22:40.899318"
$ds.Range("F85").Value = "2021-10-05 14:22:40.899321"
$ds.Range("F86").Value = "2021-10-05 14:22:40.899323"
$ds.Range("F87").Value = "2021-10-05 14:22:40.899326"
$ds.Range("F88").Value = "2021-10-05 14:22:40.899328"
$ds.Range("F89").Value = "2021-10-05 14:22:40.899331"
$ds.Range("F90").Value = "2021-10-05 14:22:40.899333"
$ds.Range("F91").Value = "2021-10-05 14:22:40.899336"
$ds.Range("F92").Value = "2021-10-05 14:22:40.899338"
$ds.Range("F93").Value = "2021-10-05 14:22:40.899341"
$ds.Range("F94").Value = "2021-10-05 14:22:40.899345"
$ds.Range("F95").Value = "2021-10-05 14:22:40.899348"
$ds.Range("F96").Value = "2021-10-05 14:22:40.899351"
$ds.Range("F97").Value = "2021-10-05 14:22:40.899353"
$ds.Range("F98").Value = "2021-10-05 14:22:40.899356"
$ds.Range("F99").Value = "2021-10-05 14:22:40.899359"
$ds.Range("F100").Value = "2021-10-05 14:22:40.899361"
$ds.Range("F101").Value = "2021-10-05 14:22:40.899364"
$ds.Range("F102").Value = "2021-10-05 14:22:40.899367"
$ds.Range("F103").Value = "2021-10-05 14:22:40.899369"
$ds.Range("F104").Value = "2021-10-05 14:22:40.899372"
$ds.Range("F105").Value = "2021-10-05 14:22:40.899374"
$ds.Range("F106").Value = "2021-10-05 14:22:40.899377"
$ds.Range("F107").Value = "2021-10-05 14:22:40.899379"
$ds.Range("F108").Value = "2021-10-05 14:22:40.899382"
$ds.Range("F109").Value = "2021-10-05 14:22:40.899385"
$ds.Range("F110").Value = "2021-10-05 14:22:40.899389"
$ds.Range("F111").Value = "2021-10-05 14:22:40.899393"
$ds.Range("F112").Value = "2021-10-05 14:22:40.899395"
$ds.Range("F113").Value = "2021-10-05 14:22:40.899398"
$ds.Range("F114").Value = "2021-10-05 14:22:40.899400"
$ds.Range("F115").Value = "2021-10-05 14:22:40.899403"
$ds.Range("F116").Value = "2021-10-05 14:22:40.899405"
$ds.Range("F117").Value = "2021-10-05 14:22:40.899408"
$ds.Range("F118").Value = "2021-10-05 14:22:40.899411"
$ds.Range("F119").Value = "2021-10-05 14:22:40.899414"
$ds.Range("F120").Value = "2021-10-05 14:22:40.899416"
$ds.Range("F121").Value = "2021-10-05 14:22:40.899419"
$ds.Range("F122").Value = "2021-10-05 14:22:40.899422"
$ds.Range("F123").Value = "2021-10-05 14:22:40.899424"
$ds.Range("F124").Value = "2021-10-05 14:22:40.899427"
$ds.Range("F125").Value = "2021-10-05 14:22:40.899429"
$ds.Range("F126").Value = "2021-10-05 14:22:40.899432"
$ds.Range("F127").Value = "2021-10-05 14:22:40.899435"
$ds.Range("F128").Value = "2021-10-05 14:22:40.899437"
$ds.Range("F129").Value = "2021-10-05 14:22:40.899440"
$ds.Range("F130").Value = "2021-10-05 14:22:40.899445"
$ds.Range("F131").Value = "2021-10-05 14:22:40.899448"
$ds.Range("F132").Value = "2021-10-05 14:22:40.899451"
$ds.Range("F133").Value = "2021-10-05 14:22:40.899453"
$ds.Range("F134").Value = "2021-10-05 14:22:40.899456"
$ds.Range("F135").Value = "2021-10-05 14:22:40.899458"
$ds.Range("F136").Value = "2021-10-05 14:22:40.899461"
$ds.Range("F137").Value = "2021-10-05 14:22:40.899464"
$ds.Range("F138").Value = "2021-10-05 14:22:40.899466"
$ds.Range("F139").Value = "2021-10-05 14:22:40.899469"
$ds.Range("F140").Value = "2021-10-05 14:22:40.899472"
$ds.Range("F141").Value = "2021-10-05 14:22:40.899474"
$ds.Range("F142").Value = "2021-10-05 14:22:40.899477"
$ds.Range("F143").Value = "2021-10-05 14:22:40.899479"
$ds.Range("F144").Value = "2021-10-05 14:22:40.899482"
$ds.Range("F145").Value = "2021-10-05 14:22:40.899485"
$ds.Range("F146").Value = "2021-10-05 14:22:40.899487"
$ds.Range("F147").Value = "2021-10-05 14:22:40.899490"
$ds.Range("F148").Value = "2021-10-05 14:22:40.899492"
$ds.Range("F149").Value = "2021-10-05 14:22:40.899495"
$ds.Range("F150").Value = "2021-10-05 14:22:40.899498"
$ds.Range("F151").Value = "2021-10-05 14:22:40.899501"
$ds.Range("F152").Value = "2021-10-05 14:22:40.899503"
$ds.Range("F153").Value = "2021-10-05 14:22:40.899506"
$ds.Range("F154").Value = "2021-10-05 14:22:40.899508"
$ds.Range("F155").Value = "2021-10-05 14:22:40.899511"
$ds.Range("F156").Value = "2021-10-05 14:22:40.899513"
$ds.Range("F157").Value = "2021-10-05 14:22:40.899516"
$ds.Range("F158").Value = "2021-10-05 14:22:40.899519"
$ds.Range("F159").Value = "2021-10-05 14:22:40.899521"
$ds.Range("F160").Value = "2021-10-05 14:22:40.899524"
$ds.Range("F161").Value = "2021-10-05 14:22:40.899527"
$ds.Range("F162").Value = "2021-10-05 14:22:40.899529"
$ds.Range("F163").Value = "2021-10-05 14:22:40.899532"
$ds.Range("F164").Value = "2021-10-05 14:22:40.899534"
$ds.Range("F165").Value = "2021-10-05 14:22:40.899537"
$ds.Range("F166").Value = "2021-10-05 14:22:40.899540"
$ds.Range("F167").Value = "2021-10-05 14:22:40.899542"
$ds.Range("F168").Value = "2021-10-05 14:22:40.899545"
$ds.Range("F169").Value = "2021-10-05 14:22:40.899547"
$ds.Range("F170").Value = "2021-10-05 14:22:40.899550"
$ds.Range("F171").Value = "2021-10-05 14:22:40.899552"
$ds.Range("F172").Value = "2021-10-05 14:22:40.899555"
$ds.Range("F173").Value = "2021-10-05 14:22:40.899557"
$ds.Range("F174").Value = "2021-10-05 14:22:40.899561"
$ds.Range("F175").Value = "2021-10-05 14:22:40.899564"
$ds.Range("F176").Value = "2021-10-05 14:22:40.899566"
$ds.Range("F177").Value = "2021-10-05 14:22:40.899569"
$ds.Range("F178").Value = "2021-10-05 14:22:40.899571"
$ds.Range("F179").Value = "2021-10-05 14:22:40.899574"
$ds.Range("F180").Value = "2021-10-05 14:22:40.899576"
$ds.Range("F181").Value = "2021-10-05 14:22:40.899579"
$ds.Range("F182").Value = "2021-10-05 14:22:40.899581"
$ds.Range("F183").Value = "2021-10-05 14:22:40.899584"
$ds.Range("F184").Value = "2021-10-05 14:22:40.899586"
$ds.Range("F185").Value = "2021-10-05 14:22:40.899588"
$ds.Range("F186").Value = "2021-10-05 14:22:40.899591"
$ds.Range("F187").Value = "2021-10-05 14:22:40.899594"
$ds.Range("F188").Value = "2021-10-05 14:22:40.899596"
$ds.Range("F189").Value = "2021-10-05 14:22:40.899598"
$ds.Range("F190").Value = "2021-10-05 14:22:40.899601"
$ds.Range("F191").Value = "2021-10-05 14:22:40.899603"
$ds.Range("F192").Value = "2021-10-05 14:22:40.899606"
$ds.Range("F193").Value = "2021-10-05 14:22:40.899608"
$ds.Range("F194").Value = "2021-10-05 14:22:40.899611"
$ds.Range("F195").Value = "2021-10-05 14:22:40.899613"
$ds.Range("F196").Value = "2021-10-05 14:22:40.899616"
$ds.Range("F197").Value = "2021-10-05 14:22:40.899618"
$ds.Range("F198").Value = "2021-10-05 14:22:40.899621"
$ds.Range("F199").Value = "2021-10-05 14:22:40.899623"
$ds.Range("F200").Value = "2021-10-05 14:22:40.899626"
$ds.Range("F201").Value = "2021-10-05 14:22:40.899628"
$ds.Range("F202").Value = "2021-10-05 14:22:40.899631"
$ds.Range("F203").Value = "2021-10-05 14:22:40.899634"
$ds.Range("F204").Value = "2021-10-05 14:22:40.899641"
$ds.Range("F205").Value = "2021-10-05 14:22:40.899643"
$ds.Range("F206").Value = "2021-10-05 14:22:40.899646"
$ds.Range("F207").Value = "2021-10-05 14:22:40.899648"
$ds.Range("F208").Value = "2021-10-05 14:22:40.899651"
$ds.Range("F209").Value = "2021-10-05 14:22:40.899653"
$ds.Range("F210").Value = "2021-10-05 14:22:40.899656"
$ds.Range("F211").Value = "2021-10-05 14:22:40.899658"
$ds.Range("F212").Value = "2021-10-05 14:22:40.899661"
$ds.Range("F213").Value = "2021-10-05 14:22:40.899663"
$ds.Range("F214").Value = "2021-10-05 14:22:40.899666"
$ds.Range("F215").Value = "2021-10-05 14:22:40.899669"
$ds.Range("F216").Value = "2021-10-05 14:22:40.899671"
$ds.Range("F217").Value = "2021-10-05 14:22:40.899674"
$ds.Range("F218").Value = "2021-10-05 14:22:40.899677"
$ds.Range("F219").Value = "2021-10-05 14:22:40.899679"
$ds.Range("F220").Value = "2021-10-05 14:22:40.899682"
$ds.Range("F221").Value = "2021-10-05 14:22:40.899684"
$ds.Range("F222").Value = "2021-10-05 14:22:40.899686"
$ds.Range("F223").Value = "2021-10-05 14:22:40.899689"
$ds.Range("F224").Value = "2021-10-05 14:22:40.899691"
$ds.Range("F225").Value = "2021-10-05 14:22:40.899694"
